$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.706.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.792.79'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.46'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.44'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.282'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0707'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.048.81'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.98'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.780.93'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.684.75'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.632'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.00'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '253.82'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0803'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.76%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.70'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.19'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.27'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.34'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.10'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0527'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.78'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.24%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.60'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.435.21'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.82%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0191'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.635'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '84.64'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.79'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.925'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.92%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.32'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.12'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.98%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.95'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.16%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.06'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0490'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.05%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.945.44'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.42'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.23%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.95'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0125'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +8.37%  '
